# Apply the "Ajout d'author dans le modele logique" edit.
#
# Summary of changes:
#  1. Metadata sheet: bump the "Date" value to the new generation timestamp.
#  2. Metadata sheet: fix the Description text - the sub-attribute list was
#     "authorInstitution, ActorXDS, authorRole ..." and becomes
#     "authorInstitution , authorPerson, authorRole ...".
#     (Elements!M2 mirrors the same Description text and updates automatically
#     because it shares the same underlying cell content.)
#  3. Elements sheet, row 5 (Author.role): cardinality relaxed from 1..1 to 0..*
#     both for the element itself (Min/Max) and for its Base Min/Base Max.
#  4. Elements sheet: the extra "Mapping: null" columns (AK, AL) - header +
#     all data cells - are removed entirely, shrinking the sheet from A1:AL6
#     to A1:AJ6.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# 1. Date
$wsMeta.Range("B8").Value = "2025-05-03T16:40:31+00:00"

# 2. Description text fix (also backs Elements!M2, same shared text)
$oldDescription = $wsMeta.Range("B12").Value()
$newDescription = $oldDescription.Replace(
    "author est un ensemble constitué des sous-attributs authorInstitution, ActorXDS, authorRole et authorSpecialty et ne porte pas de valeur par lui-même. ",
    "author est un ensemble constitué des sous-attributs authorInstitution , authorPerson, authorRole et authorSpecialty et ne porte pas de valeur par lui-même. "
)
$wsMeta.Range("B12").Value = $newDescription
$wsElem.Range("M2").Value = $newDescription

# 3. Author.role (row 5) cardinality 1..1 -> 0..*
$wsElem.Range("F5").Value = "0"
$wsElem.Range("G5").Value = "*"
$wsElem.Range("AG5").Value = "0"
$wsElem.Range("AH5").Value = "*"

# 4. Drop the "Mapping: null" columns (AK & AL) entirely.
$wsElem.Range("AK1:AL6").EntireColumn.Delete()
